# Update 2p3. Added templates for formula student suspension, torque
# vectoring, four-wheel steering.
#
# Creates a new "FSAE_Achilles" worksheet (cloned from the existing
# "Trailer_Kumanzi" template sheet, as all the Body2Axle templates share
# the same layout) and fills in the Achilles-specific numbers.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Trailer_Kumanzi")
$formatSrc = $wb.Worksheets.Item("Sedan_Hamba")

# --- clone the template sheet -------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "FSAE_Achilles"

# sheet title cell (top-right label echoing the tab name)
$ws.Range("H3").Value = "FSAE_Achilles"

# --- numeric data for the new vehicle ------------------------------------------
$ws.Range("F6").Value = -1.53

$ws.Range("F7").Value = -0.8
$ws.Range("H7").Value = 0.289
$ws.Range("K7").Value = "guesses"

$ws.Range("F8").Value = -1

$ws.Range("F9").Value = 0.25
$ws.Range("H9").Value = 0.403

$ws.Range("F10").Value = -1.75
$ws.Range("H10").Value = 0.403

$ws.Range("H11").Formula = "=0.619*2+0.2"

$ws.Range("H12").Value = 165

$ws.Range("F13").Value = 43
$ws.Range("G13").Value = 192
$ws.Range("H13").Value = 206

# --- formatting: these "guess" cells get the pink highlight used on the
# other vehicle sheets (e.g. Sedan_Hamba) instead of the template's plain
# highlight -----------------------------------------------------------------
$formatSrc.Range("F7").Copy()
$ws.Range("F7").PasteSpecial(-4122)

$formatSrc.Range("H7").Copy()
$ws.Range("H7").PasteSpecial(-4122)

$formatSrc.Range("K7").Copy()
$ws.Range("K7").PasteSpecial(-4122)

$formatSrc.Range("F8").Copy()
$ws.Range("F8").PasteSpecial(-4122)

$formatSrc.Range("F13").Copy()
$ws.Range("F13").PasteSpecial(-4122)

$formatSrc.Range("G13").Copy()
$ws.Range("G13").PasteSpecial(-4122)

$formatSrc.Range("H13").Copy()
$ws.Range("H13").PasteSpecial(-4122)

$formatSrc.Range("H11").Copy()
$ws.Range("G7").PasteSpecial(-4122)

$formatSrc.Range("H12").Copy()
$ws.Range("H12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- selections / active sheet ----------------------------------------------
$template.Range("H12").Select()

$ws.Activate()
$ws.Range("G27").Select()

$hamba = $wb.Worksheets.Item("Sedan_Hamba")
$hamba.Range("E22").Select()

$ws.Activate()
